$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $originalStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $originalStyle
}

# Enterprises density (per 1000 people) row
Set-TextValue "B11" "24.85"
Set-TextValue "C11" "3.96"
Set-TextValue "D11" "28.81"

# Employment (% of total) row
Set-TextValue "B12" "3.96"
Set-TextValue "C12" "15.78"
Set-TextValue "D12" "19.74"
